$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append after the existing last row (453), each row matching
# the no-explicit-style formatting used by the rows preceding them.
$newRows = @(
    @("www_170626b.txt",   "www_170626b.csv",   "www",    170626, "b", 12222, 6, 3, 0),
    @("wwww_170626a.txt",  "wwww_170626a.csv",  "wwww",   170626, "a", 122,   6, 3, 0),
    @("pietro_170626a.txt","pietro_170626a.csv","pietro", 170626, "a", 67,    6, 3, 0),
    @("io_170626a.txt",    "io_170626a.csv",    "io",     170626, "a", 22,    6, 3, 0),
    @("weila_170626a.txt", "weila_170626a.csv", "weila",  170626, "a", 22,    6, 3, 0)
)

$startRow = 454

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]
    for ($c = 1; $c -le 9; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $data[$c - 1]
    }
}

# The source rows already in the sheet (from row 33 onward) carry no explicit
# cell/row style (plain default formatting). Copy that "no style" formatting
# from the last pre-existing row (453) onto the freshly written rows so the
# new cells don't pick up the column's default style index.
$ws.Range("A453:I453").Copy()
$ws.Range("A454:I458").PasteSpecial(-4122)
$excel.CutCopyMode = $false
